# Applies the cryptos.xlsx price/volume/coin-name updates described by the commit diff.
# (GitHub Actions scheduled refresh of coinranking.com data, Thu Dec 21 06:27:32 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text assignment: Excel would happily store this as a string as-is
# because it is not parseable as a number (e.g. multi-dot price strings,
# percentages with surrounding spaces, coin names/links).
function Set-PlainText {
    param($addr, $val)
    $ws.Range($addr).Value = $val
}

# Text assignment for values that WOULD be auto-coerced to a number by Excel
# (e.g. "84.87"). Mark the cell as Text first so the literal string survives,
# then clear the number-format override so the cell keeps the default (General)
# style, matching how the original sheet stored these as plain text cells.
function Set-TextCell {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-PlainText 'D2' '43.620.47'
Set-PlainText 'E2' '  +1.33%  '
Set-PlainText 'D3' '2.200.61'
Set-PlainText 'E3' '  -0.89%  '
Set-PlainText 'E4' '  +0.11%  '
Set-TextCell 'D5' '258.09'
Set-PlainText 'E5' '  +1.39%  '
Set-TextCell 'D6' '84.87'
Set-PlainText 'E6' '  +11.34%  '
Set-PlainText 'E7' '  +0.68%  '
Set-PlainText 'E8' '  +0.04%  '
Set-TextCell 'D9' '0.593'
Set-PlainText 'E9' '  -0.13%  '
Set-TextCell 'D10' '44.59'
Set-PlainText 'E10' '  +7.72%  '
Set-TextCell 'D11' '0.0916'
Set-PlainText 'E11' '  -0.08%  '
Set-TextCell 'D12' '7.39'
Set-PlainText 'E12' '  +6.88%  '
Set-PlainText 'E13' '  +1.68%  '
Set-PlainText 'D14' '2.527.51'
Set-PlainText 'E14' '  -0.59%  '
Set-TextCell 'D15' '14.35'
Set-PlainText 'E15' '  -0.43%  '
Set-PlainText 'D16' '2.183.77'
Set-PlainText 'E16' '  -1.51%  '
Set-TextCell 'D17' '0.781'
Set-PlainText 'E17' '  -0.60%  '
Set-PlainText 'D18' '43.524.30'
Set-PlainText 'E18' '  +1.50%  '
Set-PlainText 'E19' '  +0.19%  '
Set-TextCell 'D20' '69.68'
Set-PlainText 'E20' '  -2.03%  '
Set-TextCell 'D21' '5.90'
Set-PlainText 'E21' '  -0.77%  '
Set-TextCell 'D22' '2.33'
Set-PlainText 'E22' '  +6.03%  '
Set-TextCell 'D23' '230.87'
Set-PlainText 'E23' '  +0.54%  '
Set-TextCell 'D24' '8.92'
Set-PlainText 'E24' '  -4.80%  '
Set-PlainText 'E25' '  -0.02%  '
Set-TextCell 'D26' '3.60'
Set-PlainText 'E26' '  +6.78%  '
Set-TextCell 'D27' '10.61'
Set-PlainText 'E27' '  +0.04%  '
Set-PlainText 'E28' '  +2.81%  '
Set-PlainText 'B29' 'InjectiveProtocol'
Set-PlainText 'C29' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 'D29' '39.21'
Set-PlainText 'E29' '  +1.52%  '
Set-PlainText 'B30' 'PancakeSwap'
Set-PlainText 'C30' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D30' '2.23'
Set-PlainText 'E30' '  +1.95%  '
Set-TextCell 'D31' '173.44'
Set-PlainText 'E31' '  +0.17%  '
Set-TextCell 'D32' '20.36'
Set-PlainText 'E32' '  +0.46%  '
Set-TextCell 'D33' '0.0865'
Set-PlainText 'E33' '  +2.48%  '
Set-TextCell 'D34' '5.32'
Set-PlainText 'E34' '  +1.67%  '
Set-PlainText 'E35' '  +1.23%  '
Set-TextCell 'D36' '0.110'
Set-PlainText 'E36' '  +0.39%  '
Set-PlainText 'B37' 'RenderToken'
Set-PlainText 'C37' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D37' '4.47'
Set-PlainText 'E37' '  +3.82%  '
Set-PlainText 'B38' 'VeChain'
Set-PlainText 'C38' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D38' '0.0356'
Set-PlainText 'E38' '  +1.48%  '
Set-PlainText 'B39' 'Celestia'
Set-PlainText 'C39' 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextCell 'D39' '12.43'
Set-PlainText 'E39' '  -0.70%  '
Set-PlainText 'B40' 'NEARProtocol'
Set-PlainText 'C40' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D40' '2.86'
Set-PlainText 'E40' '  +4.50%  '
Set-TextCell 'D41' '2.09'
Set-PlainText 'E41' '  -0.90%  '
Set-TextCell 'D42' '63.14'
Set-PlainText 'E42' '  +5.02%  '
Set-TextCell 'D43' '5.44'
Set-PlainText 'E43' '  +3.42%  '
Set-TextCell 'D44' '0.198'
Set-PlainText 'E44' '  -0.11%  '
Set-PlainText 'B45' 'FraxShare'
Set-PlainText 'C45' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D45' '8.35'
Set-PlainText 'E45' '  +0.18%  '
Set-TextCell 'D46' '100.05'
Set-PlainText 'E46' '  -1.85%  '
Set-PlainText 'B47' 'Cronos'
Set-PlainText 'C47' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D47' '0.0979'
Set-PlainText 'E47' '  +0.13%  '
Set-TextCell 'D48' '1.18'
Set-PlainText 'E48' '  +3.61%  '
Set-PlainText 'E49' '  -0.13%  '
Set-TextCell 'D50' '0.437'
Set-PlainText 'E50' '  -5.27%  '
Set-PlainText 'E51' '  +2.05%  '
